# This edit reorders the weekly price records (rows 2-12) in the sheet.
# The per-row "dimension" columns (A,B,C,E,F,G,H,I,J) are identical for every
# record, so only the varying columns (D,K,L,M,N,O,P,Q,R,S,T) need to move
# between rows. We snapshot the current values first, then write them back
# out according to the new row order described by the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("D", "K", "L", "M", "N", "O", "P", "Q", "R", "S", "T")

# Take a snapshot of the current (pre-edit) values for the varying columns
# of every data row (2 through 12).
$snapshot = @{}
for ($r = 2; $r -le 12; $r++) {
    $rowData = @{}
    foreach ($col in $cols) {
        $rowData[$col] = $ws.Range("$col$r").Value2
    }
    $snapshot[$r] = $rowData
}

# Mapping of new row number -> original row number that its data came from.
$mapping = @{
    2  = 12
    3  = 5
    4  = 4
    5  = 3
    6  = 8
    7  = 9
    8  = 10
    9  = 6
    10 = 11
    11 = 7
    12 = 2
}

foreach ($newRow in $mapping.Keys) {
    $srcRow = $mapping[$newRow]
    $srcData = $snapshot[$srcRow]
    foreach ($col in $cols) {
        $ws.Range("$col$newRow").Value = $srcData[$col]
    }
}
